$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.703.28"
$ws.Range("E2").Value = "  -1.47%  "
$ws.Range("D3").Value = "3.170.22"
$ws.Range("E3").Value = "  -3.62%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.57"
$ws.Range("E5").Value = "  -1.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.37"
$ws.Range("E6").Value = "  -3.85%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.167.02"
$ws.Range("E8").Value = "  -3.73%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.509"
$ws.Range("E9").Value = "  -1.99%  "
$ws.Range("E10").Value = "  -4.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.28"
$ws.Range("E11").Value = "  -3.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.456"
$ws.Range("E12").Value = "  -3.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000235"
$ws.Range("E13").Value = "  -4.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.85"
$ws.Range("E14").Value = "  -1.89%  "
$ws.Range("D15").Value = "3.695.51"
$ws.Range("E15").Value = "  -3.54%  "
$ws.Range("E16").Value = "  -2.09%  "
$ws.Range("D17").Value = "3.171.52"
$ws.Range("E17").Value = "  -3.53%  "
$ws.Range("D18").Value = "62.661.73"
$ws.Range("E18").Value = "  -1.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.57"
$ws.Range("E19").Value = "  -3.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "455.35"
$ws.Range("E20").Value = "  -4.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.00"
$ws.Range("E21").Value = "  -0.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.703"
$ws.Range("E22").Value = "  -3.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.63"
$ws.Range("E23").Value = "  -5.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.42"
$ws.Range("E24").Value = "  -1.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.77"
$ws.Range("E25").Value = "  -0.59%  "
$ws.Range("E27").Value = "  -2.31%  "
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.76"
$ws.Range("E29").Value = "  -3.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.78"
$ws.Range("E30").Value = "  -6.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.03"
$ws.Range("E31").Value = "  -5.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.37"
$ws.Range("E32").Value = "  -4.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.103"
$ws.Range("E33").Value = "  -2.14%  "
$ws.Range("E34").Value = "  -5.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.85"
$ws.Range("E36").Value = "  -2.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.23"
$ws.Range("E37").Value = "  -4.00%  "
$ws.Range("D38").Value = "0.0₃0709"
$ws.Range("E38").Value = "  -3.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0387"
$ws.Range("E39").Value = "  -2.73%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "406.16"
$ws.Range("E40").Value = "  -5.48%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.71"
$ws.Range("E41").Value = "  -0.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.05"
$ws.Range("E42").Value = "  -3.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.112"
$ws.Range("E43").Value = "  -1.03%  "
$ws.Range("D44").Value = "2.797.71"
$ws.Range("E44").Value = "  -8.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.251"
$ws.Range("E45").Value = "  -5.03%  "
$ws.Range("E46").Value = "  -1.97%  "
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "126.25"
$ws.Range("E48").Value = "  -1.17%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.50"
$ws.Range("E49").Value = "  -2.74%  "
$ws.Range("B50").Value = "Arweave"
$ws.Range("C50").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.89"
$ws.Range("E50").Value = "  -1.90%  "
$ws.Range("E51").Value = "  -3.04%  "
